$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (it will no longer exist in the final data range A1:E10)
$ws.Rows.Item(11).Delete()

# Vendor names used for rows 2-10 (cycle of 3)
$vendors = @("Nestle Sofia Corp.", "Zagorka Corp.", "Targovishte Bottling Company Ltd.")

# Data values for incomes, expenses, total taxes, financial result per vendor pattern
$data = @{
    "Nestle Sofia Corp." = @(200, 30, 36, 134)
    "Zagorka Corp." = @(100, 120, 20, -40)
    "Targovishte Bottling Company Ltd." = @(100, 200, 25, -125)
}

for ($i = 0; $i -lt 9; $i++) {
    $row = 2 + $i
    $vendor = $vendors[$i % 3]
    $vals = $data[$vendor]

    $ws.Cells.Item($row, 1).Value = $vendor
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
